$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (NN model results) for rows 2-9, columns B-F and H.
# Column G (Площадь/Area) is unchanged.
$data = @{
    2 = @{ B = 1308.811279296875;  C = 0.9819;  D = 0.932699978351593;  E = 1.58050000667572;   F = 0.7056000232696533; H = 1.0029 }
    3 = @{ B = 1202.23046875;      C = 0.9618;  D = 0.9407;             E = 2.047100067138672;   F = 0.7408999800682068; H = 1.0737 }
    4 = @{ B = 730.177490234375;   C = 0.853;   D = 0.8435;             E = 1.135800004005432;   F = 0.7342000007629395; H = 0.2126 }
    5 = @{ B = 727.30419921875;    C = 0.7672;  D = 0.761;              E = 1.112400054931641;   F = 0.6722999811172485; H = -0.5187 }
    6 = @{ B = 1057.1865234375;    C = 0.8364;  D = 0.8395;             E = 1.127099990844727;   F = 0.7050999999046326; H = 0.1768 }
    7 = @{ B = 913.6450805664062;  C = 0.9201;  D = 0.9129999876022339; E = 1.260900020599365;   F = 0.7613999843597412; H = 0.8276 }
    8 = @{ B = 1033.1396484375;    C = 0.9258;  D = 0.9162;             E = 1.295899987220764;   F = 0.8040000200271606; H = 0.8563 }
    9 = @{ B = 6972.49462890625;   C = 0.8986;  D = 0.8848;             E = 2.047100067138672;   F = 0.6722999811172485; H = 3.6312 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}
